$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Add new row 10 with the download section description (register this shared string first)
$ws.Cells.Item(10, 2).Value = "Downloadsectie selecties aanmaken met tekst over wat ze inhouden en wat je download"

# Fix capitalization of the "uitwerken/invoegen..." row (row 9, column B)
$ws.Cells.Item(9, 2).Value = "Uitwerken/invoegen algemene info (aangeleverd in een word bestand)"

# Update selection to match target
$ws.Range("B18").Select()
